# Bug Fix - CellTempMaxActual
#
# The sheet had four bogus "CV_DTC1".."CV_DTC4" header columns (ET:EW)
# inserted ahead of the real CellVoltageMinActual/CellVoltageMaxActual/
# CellTempMinActual/CellVoltageMaxActual (CellTempMaxActual) block, which
# pushed every later column one set too far right and left the real
# values for that block unpopulated. Remove the bogus columns (shifting
# everything after them back into place) and fill in the now-correctly
# positioned CellVoltageMinActual / CellVoltageMaxActual / CellTempMinActual
# / CellTempMaxActual / TempVoltageResult values for row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the 4 spurious CV_DTC1..CV_DTC4 columns (ET:EW); everything to
# the right shifts left by 4, which is exactly what the diff shows
# (dimension A1:GA2 -> A1:FW2, and every header/value past column ES
# moving back four slots).
$ws.Range("ET1:EW1").EntireColumn.Delete()

# These six cells land on the now-shifted ES2:EX2 slots and previously
# held no data; populate them with the actual measured values. Force
# text formatting first so the numeric-looking strings ("4.0", "3.631",
# ...) are stored as text instead of being auto-converted to numbers,
# matching the inline-string cells used throughout this sheet.
$ws.Range("ES2:EX2").NumberFormat = "@"
$ws.Range("ES2").Value = "4.0"
$ws.Range("ET2").Value = "3.631"
$ws.Range("EU2").Value = "3.635"
$ws.Range("EV2").Value = "22.500"
$ws.Range("EW2").Value = "3.635"
$ws.Range("EX2").Value = '"IO"'
